# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy formatting from an existing header cell (AB1) so the
# new header cells (AC1, AD1, AE1) pick up the same bold/centered/bordered style.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value2 = "Wins"
$ws.Range("AD1").Value2 = "Losses"
$ws.Range("AE1").Value2 = "Ties"

# Data rows 2-41: every row gets the same season record values.
$wins = 93
$losses = 69
$ties = 0

for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 29).Value2 = $wins    # column AC
    $ws.Cells.Item($r, 30).Value2 = $losses  # column AD
    $ws.Cells.Item($r, 31).Value2 = $ties    # column AE
}
